$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 478.375
$ws.Range("I41").Value = 225
$ws.Range("K41").Value = 225
$ws.Range("M41").Value = 215
# Row 86
$ws.Range("H86").Value = 4694
$ws.Range("I86").Value = 2496.75
$ws.Range("J86").Value = 6891.25
$ws.Range("K86").Value = 2496.75
$ws.Range("L86").Value = 6891.25
$ws.Range("M86").Value = -1373.75
$ws.Range("N86").Value = -9137.25
# Row 89
$ws.Range("H89").Value = 4694
$ws.Range("I89").Value = 2496.75
$ws.Range("J89").Value = 6891.25
$ws.Range("K89").Value = 12483.75
$ws.Range("L89").Value = 34456.25
$ws.Range("M89").Value = -6867.75
$ws.Range("N89").Value = -45688.25
# Row 100
$ws.Range("H100").Value = 724.6667
$ws.Range("I100").Value = 724.6667
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 724.6667
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -183.6667
$ws.Range("N100").ClearContents()
# Row 131
$ws.Range("H131").Value = 2363
$ws.Range("I131").Value = 2363
$ws.Range("K131").Value = 7089
$ws.Range("M131").Value = -2049
# Row 132
$ws.Range("H132").Value = 7599.4
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
# Row 137
$ws.Range("H137").Value = 1104.1765
$ws.Range("J137").Value = 1265.1111
$ws.Range("L137").Value = 3795.3333
$ws.Range("N137").Value = -8895.3333
# Row 141
$ws.Range("H141").Value = 2985.2856
$ws.Range("I141").Value = 3059.8
$ws.Range("J141").Value = 2799
$ws.Range("K141").Value = 9179.400000000001
$ws.Range("L141").Value = 8397
$ws.Range("M141").Value = -3999.400000000001
$ws.Range("N141").Value = -18757

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 918.4
$ws.Range("I61").Value = 916.75
$ws.Range("K61").Value = 916.75
$ws.Range("M61").Value = -704.75
# Row 63
$ws.Range("H63").Value = 1290
$ws.Range("I63").Value = 1290
$ws.Range("K63").Value = 1290
$ws.Range("M63").Value = -604
# Row 66
$ws.Range("H66").Value = 1290
$ws.Range("I66").Value = 1290
$ws.Range("K66").Value = 6450
$ws.Range("M66").Value = -3018
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# Row 136
$ws.Range("H136").Value = 918.4
$ws.Range("I136").Value = 916.75
$ws.Range("K136").Value = 2750.25
$ws.Range("M136").Value = -200.25

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2899.8
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
# Row 107
$ws.Range("H107").Value = 196
$ws.Range("I107").Value = 196
$ws.Range("K107").Value = 196
$ws.Range("M107").Value = 1724
# Row 134
$ws.Range("H134").Value = 4926.857
$ws.Range("I134").Value = 4882.769
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 14648.307
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -12113.307
$ws.Range("N134").Value = -21570

$ws = $wb.Worksheets.Item("CRP")
# Row 92
$ws.Range("H92").Value = 26916
$ws.Range("J92").Value = 26916
$ws.Range("L92").Value = 26916
$ws.Range("N92").Value = -31908
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# Row 134
$ws.Range("H134").Value = 2185.5
$ws.Range("I134").Value = 2366.6667
$ws.Range("J134").Value = 2004.3334
$ws.Range("K134").Value = 7100.000100000001
$ws.Range("L134").Value = 6013.0002
$ws.Range("M134").Value = -4565.000100000001
$ws.Range("N134").Value = -11083.0002

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1156
$ws.Range("J5").Value = 1092.7142
$ws.Range("L5").Value = 3278.1426
$ws.Range("N5").Value = -3502.1426
# Row 81
$ws.Range("H81").Value = 2500
$ws.Range("I81").Value = 2500
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 7500
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -6377
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 2500
$ws.Range("I84").Value = 2500
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 22500
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -16884
$ws.Range("N84").ClearContents()
# Row 121
$ws.Range("H121").Value = 353.33334
$ws.Range("I121").Value = 353.33334
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 1060.00002
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 249.9999800000001
$ws.Range("N121").ClearContents()
# Row 135
$ws.Range("H135").Value = 1156
$ws.Range("J135").Value = 1092.7142
$ws.Range("L135").Value = 9834.427799999999
$ws.Range("N135").Value = -14904.4278
# Row 139
$ws.Range("H139").Value = 6880
$ws.Range("I139").Value = 6880
$ws.Range("K139").Value = 20640
$ws.Range("M139").Value = -15500
# Row 140
$ws.Range("H140").Value = 1035.2222
$ws.Range("I140").Value = 1035.2222
$ws.Range("K140").Value = 3105.6666
$ws.Range("M140").Value = 2074.3334

$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 7505000
$ws.Range("I35").Value = 7505000
$ws.Range("K35").Value = 7505000
$ws.Range("M35").Value = -7504702
# Row 70
$ws.Range("H70").Value = 4676.2
$ws.Range("I70").Value = 4643.6665
$ws.Range("K70").Value = 4643.6665
$ws.Range("M70").Value = -4373.6665
# Row 73
$ws.Range("H73").Value = 4676.2
$ws.Range("I73").Value = 4643.6665
$ws.Range("K73").Value = 4643.6665
$ws.Range("M73").Value = -3707.6665
# Row 97
$ws.Range("H97").Value = 2912
$ws.Range("I97").Value = 2447.125
$ws.Range("K97").Value = 2447.125
$ws.Range("M97").Value = -1951.125
# Row 132
$ws.Range("H132").Value = 3862.75
$ws.Range("I132").Value = 3680.138
$ws.Range("K132").Value = 11040.414
$ws.Range("M132").Value = -8510.414000000001
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 32284.572
$ws.Range("I2").Value = 8666.666999999999
$ws.Range("K2").Value = 8666.666999999999
$ws.Range("M2").Value = -8554.666999999999
# Row 22
$ws.Range("H22").Value = 1274.25
$ws.Range("I22").Value = 838.8
$ws.Range("K22").Value = 838.8
$ws.Range("M22").Value = -543.8
# Row 27
$ws.Range("H27").Value = 1274.25
$ws.Range("I27").Value = 838.8
$ws.Range("K27").Value = 838.8
$ws.Range("M27").Value = -731.8
# Row 33
$ws.Range("H33").Value = 62946.668
$ws.Range("J33").Value = 69420
$ws.Range("L33").Value = 69420
$ws.Range("N33").Value = -70000
# Row 40
$ws.Range("H40").Value = 4245.5
$ws.Range("I40").Value = 4245.5
$ws.Range("K40").Value = 4245.5
$ws.Range("M40").Value = -4109.5
# Row 46
$ws.Range("H46").Value = 3849.9412
$ws.Range("I46").Value = 3313.5454
$ws.Range("K46").Value = 3313.5454
$ws.Range("M46").Value = -3125.5454
# Row 100
$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 1500
$ws.Range("M100").Value = -959
# Row 122
$ws.Range("H122").Value = 4633
$ws.Range("I122").Value = 4633
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13899
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11449
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 5000225
$ws.Range("I20").Value = 450
$ws.Range("K20").Value = 450
$ws.Range("M20").Value = -210
# Row 100
$ws.Range("H100").Value = 3669554.5
$ws.Range("I100").Value = 4979273
$ws.Range("J100").Value = 2343
$ws.Range("K100").Value = 9958546
$ws.Range("L100").Value = 4686
$ws.Range("M100").Value = -9958005
$ws.Range("N100").Value = -5768
# Row 126
$ws.Range("H126").Value = 912.8889
$ws.Range("I126").Value = 1002.61536
$ws.Range("J126").Value = 679.6
$ws.Range("K126").Value = 3007.84608
$ws.Range("L126").Value = 2038.8
$ws.Range("M126").Value = -537.8460800000003
$ws.Range("N126").Value = -6978.8
# Row 132
$ws.Range("H132").Value = 5500
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
